$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 8
$ws.Range("F4").Value = 1
$ws.Range("F7").Value = -7
$ws.Range("F9").Value = -3
$ws.Range("F11").Value = 3
